# Updated CVDs for the month
# Clears / updates the Jul ("O") Commit/Forecast figures (CVD rolled off) and
# recomputes the downstream Aug/Sep/Q3/Oct/Nov/Dec/Q4/FY figures for the
# affected (location, cvd, data_source="Commit/Forecast") rows.

$wb = $excel.ActiveWorkbook

# --- Sheet "Baja California Mexico" : Professional Voluntary Turnover, row 4 ---
$ws = $wb.Worksheets.Item("Baja California Mexico")
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

# --- Sheet "Baja California Mexico" : Manufacturing Voluntary Turnover, row 9 ---
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("O9").ClearContents()

# --- Sheet "Charlotte  North Carolina" : row 4 ---
$ws = $wb.Worksheets.Item("Charlotte  North Carolina")
$ws.Range("O4").ClearContents()

# --- Sheet "Cleveland Ohio" : row 6 ---
$ws = $wb.Worksheets.Item("Cleveland Ohio")
$ws.Range("O6").ClearContents()

# --- Sheet "Marengo Illinois" : row 4 and row 9 ---
$ws = $wb.Worksheets.Item("Marengo Illinois")
$ws.Range("O4").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("O9").ClearContents()
